# Rename sheets:
#   Hoja1 -> Migrante
#   Hoja2 -> arg1
#   Hoja3 -> arg2
#   Hoja4 -> total
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Hoja1").Name = "Migrante"
$wb.Worksheets.Item("Hoja2").Name = "arg1"
$wb.Worksheets.Item("Hoja3").Name = "arg2"
$wb.Worksheets.Item("Hoja4").Name = "total"

# Scroll the "Output" sheet so row 6 is at the top (was row 22), keep the
# existing C30:F34 selection.
$wsOutput = $wb.Worksheets.Item("Output")
$wsOutput.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$wsOutput.Range("C30:F34").Select()

# The previously active sheet was "arg1" (was "Hoja2"); the new active /
# selected sheet is "total" (was "Hoja4").
$wb.Worksheets.Item("total").Activate()
